$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date shifts from 44692 -> 44691, volume 120 -> 100
$ws.Range("D2").Value = 44691
$ws.Range("J2").Value = 100

# Row 3: date shifts from 44691 -> 44687, volume 100 -> 160
$ws.Range("D3").Value = 44687
$ws.Range("J3").Value = 160

# Row 4: date shifts from 44687 -> 44221, and full price/origin block
# changes to the "atado" / Provincia de Diguillin record
$ws.Range("D4").Value = 44221
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 1300
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = 1420
$ws.Range("N4").Value = "$/atado"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 1420
$ws.Range("Q4").Value = 1

# Row 5: date shifts from 44221 -> 44692, and full price/origin block
# changes to the "docena de matas" / Region Metropolitana record
$ws.Range("D5").Value = 44692
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3500
$ws.Range("M5").Value = 3250
$ws.Range("N5").Value = "$/docena de matas"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 542
$ws.Range("Q5").Value = 6
